# "Started adding async methods"
# Adds new checklist columns (Synchronus / Asynchronus / .Net 3.5 / Integration
# Test / Documentation) to the "Methods" sheet, and marks a handful of rows
# as "Done" (using the existing green "Good" cell style) in columns A, C, D
# for the entries that already had progress tracked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

# New header row for the added tracking columns.
$ws.Range("C1").Value = "Synchronus"
$ws.Range("D1").Value = "Asynchronus"
$ws.Range("E1").Value = ".Net 3.5"
$ws.Range("F1").Value = "Integration Test"
$ws.Range("G1").Value = "Documentation"

# Rows whose "Entity" cell (column A) needs to be (newly) marked Done.
$doneEntityRows = @(2, 3, 88)
foreach ($r in $doneEntityRows) {
    $cell = $ws.Range("A$r")
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Rows whose Synchronus/Asynchronus columns (C/D) are marked Done.
$doneCDRows = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 13, 14, 81, 82, 83, 84)
foreach ($r in $doneCDRows) {
    $ws.Range("C$r").Value = "Done"
    $ws.Range("D$r").Value = "Done"
}

# New columns sized (bestFit-equivalent) to fit their header/content.
$ws.Columns.Item(3).ColumnWidth = 10.25
$ws.Columns.Item(4).ColumnWidth = 11.416666666666666
$ws.Columns.Item(5).ColumnWidth = 6.916666666666667
$ws.Columns.Item(6).ColumnWidth = 14.25
$ws.Columns.Item(7).ColumnWidth = 13.916666666666666

# Restore the view state (selected cell) recorded in the edited workbook.
$ws.Activate()
$ws.Range("D81").Select()
